$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in columns E and F
$ws.Range("E2").Value = 0.58899999999999997
$ws.Range("F2").Value = 0.57699999999999996

$ws.Range("E3").Value = 8.9670000000000005

$ws.Range("E4").Value = 8.4920000000000009

$ws.Range("E5").Value = 1.4319999999999999
$ws.Range("F5").Value = 0.29399999999999998

$ws.Range("E6").Value = 7.7889999999999997

$ws.Range("E7").Value = 0.40799999999999997
$ws.Range("F7").Value = 0.86699999999999999

$ws.Range("E8").Value = 0.33100000000000002
$ws.Range("F8").Value = 0.72799999999999998

$ws.Range("E9").Value = 5.5839999999999996
$ws.Range("F9").Value = 0.005

$ws.Range("E10").Value = 12.499000000000001

# F9 becomes bold (new style)
$ws.Range("F9").Font.Bold = $true

# Update selection to I7
$ws.Range("I7").Select()
